$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the data range so numeric-looking strings
# (e.g. "0.9980", "1.262.21") are preserved exactly as typed,
# matching the source workbook where these are stored as strings.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '29.348.57'
$ws.Range('E2').Value = '  +0.02%  '

$ws.Range('D3').Value = '1.844.83'
$ws.Range('E3').Value = '  -0.10%  '

$ws.Range('D4').Value = '0.9980'
$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').Value = '239.94'
$ws.Range('E5').Value = '  -0.08%  '

$ws.Range('E6').Value = '  +0.15%  '

$ws.Range('D7').Value = '0.9995'
$ws.Range('E7').Value = '  +0.07%  '

$ws.Range('D8').Value = '0.07488'

$ws.Range('E9').Value = '  -0.13%  '

$ws.Range('D10').Value = '24.49'
$ws.Range('E10').Value = '  -0.88%  '

$ws.Range('D11').Value = '0.07737'
$ws.Range('E11').Value = '  +0.04%  '

$ws.Range('D12').Value = '1.845.33'
$ws.Range('E12').Value = '  -2.27%  '

$ws.Range('D13').Value = '4.985'
$ws.Range('E13').Value = '  -0.78%  '

$ws.Range('E14').Value = '  +0.18%  '

$ws.Range('D15').Value = '0.00001052'
$ws.Range('E15').Value = '  -0.69%  '

$ws.Range('D16').Value = '81.99'
$ws.Range('E16').Value = '  -1.14%  '

$ws.Range('D17').Value = '6.169'
$ws.Range('E17').Value = '  +0.77%  '

$ws.Range('D18').Value = '29.377.34'
$ws.Range('E18').Value = '  -0.02%  '

$ws.Range('D19').Value = '229.07'
$ws.Range('E19').Value = '  +0.62%  '

$ws.Range('E20').Value = '  -0.08%  '

$ws.Range('D21').Value = '0.9990'
$ws.Range('E21').Value = '  +0.04%  '

$ws.Range('D22').Value = '7.503'
$ws.Range('E22').Value = '  +0.45%  '

$ws.Range('D23').Value = '0.9991'
$ws.Range('E23').Value = '  +0.00%  '

$ws.Range('D24').Value = '158.50'
$ws.Range('E24').Value = '  -0.13%  '

$ws.Range('E26').Value = '  -0.85%  '

$ws.Range('D27').Value = '17.53'
$ws.Range('E27').Value = '  -0.63%  '

$ws.Range('D28').Value = '0.06527'
$ws.Range('E28').Value = '  +16.68%  '

$ws.Range('D29').Value = '1.411'
$ws.Range('E29').Value = '  -1.65%  '

$ws.Range('E30').Value = '  +1.37%  '

$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').Value = '4.109'
$ws.Range('E31').Value = '  +1.07%  '

$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '4.091'
$ws.Range('E32').Value = '  -0.28%  '

$ws.Range('E33').Value = '  -0.03%  '

$ws.Range('E34').Value = '  -1.67%  '

$ws.Range('D35').Value = '0.6964'
$ws.Range('E35').Value = '  +0.14%  '

$ws.Range('E36').Value = '  -0.09%  '

$ws.Range('D37').Value = '1.262.21'
$ws.Range('E37').Value = '  +2.96%  '

$ws.Range('D38').Value = '2.835'
$ws.Range('E38').Value = '  +4.32%  '

$ws.Range('E39').Value = '  +2.01%  '

$ws.Range('D40').Value = '6.766'
$ws.Range('E40').Value = '  +6.51%  '

$ws.Range('D41').Value = '0.9188'
$ws.Range('E41').Value = '  +2.54%  '

$ws.Range('D42').Value = '0.9986'
$ws.Range('E42').Value = '  +0.01%  '

$ws.Range('D43').Value = '2.007.78'
$ws.Range('E43').Value = '  +1.47%  '

$ws.Range('D44').Value = '101.26'
$ws.Range('E44').Value = '  -0.12%  '

$ws.Range('D45').Value = '66.16'
$ws.Range('E45').Value = '  +1.08%  '

$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '0.00000000119'
$ws.Range('E46').Value = '  +0.12%  '

$ws.Range('D47').Value = '1.726'
$ws.Range('E47').Value = '  +2.23%  '

$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').Value = '7.078'
$ws.Range('E48').Value = '  -1.83%  '

$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value = '0.1162'
$ws.Range('E49').Value = '  +1.95%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '8.984'
$ws.Range('E50').Value = '  -0.45%  '

$ws.Range('B51').Value = 'TheSandbox'
$ws.Range('C51').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D51').Value = '0.3955'
$ws.Range('E51').Value = '  -0.73%  '

# Restore default "Normal" style so cells have no explicit
# number-format override left behind (matches source formatting).
$dataRange.Style = "Normal"
